{"js": "// Update the worksheet's \"two-digit \u00f7 one-digit\" problems to a newly\n// generated set of numbers. The document has a single 20-row x 5-column\n// table where only every 4th row (0, 4, 8, 12, 16) holds a division\n// problem such as \"46\u00f73=\" \u2014 the rows in between are left blank for the\n// student's work. Each populated cell's text is replaced in place while\n// keeping the existing run/paragraph formatting (fonts, size, alignment)\n// untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, columnIndex, newText] \u2014 order follows the table's reading\n// order (row-major), matching the order of replacements in the diff.\nconst replacements = [\n  [0, 0, \"39\u00f73=\"],\n  [0, 1, \"48\u00f75=\"],\n  [0, 2, \"41\u00f75=\"],\n  [0, 3, \"30\u00f76=\"],\n  [0, 4, \"80\u00f76=\"],\n  [4, 0, \"42\u00f74=\"],\n  [4, 1, \"55\u00f76=\"],\n  [4, 2, \"20\u00f77=\"],\n  [4, 3, \"47\u00f72=\"],\n  [4, 4, \"82\u00f78=\"],\n  [8, 0, \"90\u00f79=\"],\n  [8, 1, \"42\u00f72=\"],\n  [8, 2, \"53\u00f78=\"],\n  [8, 3, \"55\u00f74=\"],\n  [8, 4, \"22\u00f77=\"],\n  [12, 0, \"83\u00f76=\"],\n  [12, 1, \"46\u00f77=\"],\n  [12, 2, \"75\u00f76=\"],\n  [12, 3, \"42\u00f73=\"],\n  [12, 4, \"24\u00f78=\"],\n  [16, 0, \"42\u00f75=\"],\n  [16, 1, \"44\u00f79=\"],\n  [16, 2, \"86\u00f72=\"],\n  [16, 3, \"33\u00f75=\"],\n  [16, 4, \"26\u00f77=\"],\n];\n\n// Grab the first paragraph of each target cell up front, then mutate.\nconst paragraphs = [];\nfor (const [row, col] of replacements) {\n  const cell = table.getCell(row, col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraphs.push(paragraph);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const newText = replacements[i][2];\n  const range = paragraphs[i].getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the worksheet's \"two-digit \u00f7 one-digit\" problems to a newly\n# generated set of numbers. The document has a single 20-row x 5-column\n# table where only every 4th row (1, 5, 9, 13, 17 in 1-based COM terms)\n# holds a division problem such as \"46\u00f73=\" -- the rows in between are\n# left blank for the student's work. Each populated cell's text is\n# replaced in place (via Cell.Range.Text) so the existing run/paragraph\n# formatting (fonts, size, alignment) is left untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# row, column (1-based), newText -- row-major order, matching the diff.\n$replacements = @(\n    @(1, 1, \"39\u00f73=\"),\n    @(1, 2, \"48\u00f75=\"),\n    @(1, 3, \"41\u00f75=\"),\n    @(1, 4, \"30\u00f76=\"),\n    @(1, 5, \"80\u00f76=\"),\n    @(5, 1, \"42\u00f74=\"),\n    @(5, 2, \"55\u00f76=\"),\n    @(5, 3, \"20\u00f77=\"),\n    @(5, 4, \"47\u00f72=\"),\n    @(5, 5, \"82\u00f78=\"),\n    @(9, 1, \"90\u00f79=\"),\n    @(9, 2, \"42\u00f72=\"),\n    @(9, 3, \"53\u00f78=\"),\n    @(9, 4, \"55\u00f74=\"),\n    @(9, 5, \"22\u00f77=\"),\n    @(13, 1, \"83\u00f76=\"),\n    @(13, 2, \"46\u00f77=\"),\n    @(13, 3, \"75\u00f76=\"),\n    @(13, 4, \"42\u00f73=\"),\n    @(13, 5, \"24\u00f78=\"),\n    @(17, 1, \"42\u00f75=\"),\n    @(17, 2, \"44\u00f79=\"),\n    @(17, 3, \"86\u00f72=\"),\n    @(17, 4, \"33\u00f75=\"),\n    @(17, 5, \"26\u00f77=\")\n)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $newText = $entry[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $newText\n}\n"}
